$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 303.1111
$ws.Range("I2").Value = 228.14285
$ws.Range("J2").Value = 565.5
$ws.Range("K2").Value = 228.14285
$ws.Range("L2").Value = 565.5
$ws.Range("M2").Value = -115.14285
$ws.Range("N2").Value = -791.5
$ws.Range("H39").Value = 474.10526
$ws.Range("I39").Value = 35.57143
$ws.Range("J39").Value = 729.9167
$ws.Range("K39").Value = 106.71429
$ws.Range("L39").Value = 2189.7501
$ws.Range("M39").Value = 189.28571
$ws.Range("N39").Value = -2781.7501
$ws.Range("H107").Value = 624.2105
$ws.Range("I107").Value = 656.05
$ws.Range("J107").Value = 549.2941
$ws.Range("K107").Value = 656.05
$ws.Range("L107").Value = 549.2941
$ws.Range("M107").Value = 1263.95
$ws.Range("N107").Value = -4389.2941
$ws.Range("H125").Value = 2600.6667
$ws.Range("J125").Value = 2920.8
$ws.Range("L125").Value = 26287.2
$ws.Range("N125").Value = -31207.2

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 182.83333
$ws.Range("I5").Value = 224.25
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 224.25
$ws.Range("L5").Value = 100
$ws.Range("M5").Value = -112.25
$ws.Range("N5").Value = -324
$ws.Range("H32").Value = 14711021
$ws.Range("I32").Value = 20409868
$ws.Range("J32").Value = 13995
$ws.Range("K32").Value = 20409868
$ws.Range("L32").Value = 13995
$ws.Range("M32").Value = -20409581
$ws.Range("N32").Value = -14569
$ws.Range("H122").Value = 888
$ws.Range("I122").Value = 791.7692
$ws.Range("J122").Value = 1513.5
$ws.Range("K122").Value = 2375.3076
$ws.Range("L122").Value = 4540.5
$ws.Range("M122").Value = 74.69239999999991
$ws.Range("N122").Value = -9440.5
$ws.Range("H132").Value = 3354.4194
$ws.Range("I132").Value = 3047.9092
$ws.Range("J132").Value = 4103.6665
$ws.Range("K132").Value = 9143.7276
$ws.Range("L132").Value = 12310.9995
$ws.Range("M132").Value = -6613.7276
$ws.Range("N132").Value = -17370.9995

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 182.83333
$ws.Range("I4").Value = 224.25
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 224.25
$ws.Range("L4").Value = 100
$ws.Range("M4").Value = -109.25
$ws.Range("N4").Value = -330
$ws.Range("H44").Value = 13266.667
$ws.Range("I44").Value = 10000
$ws.Range("J44").Value = 14900
$ws.Range("K44").Value = 10000
$ws.Range("L44").Value = 14900
$ws.Range("M44").Value = -9503
$ws.Range("N44").Value = -15894

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 49.333332
$ws.Range("I7").Value = 48.88889
$ws.Range("J7").Value = 49.77778
$ws.Range("K7").Value = 48.88889
$ws.Range("L7").Value = 49.77778
$ws.Range("M7").Value = 64.11111
$ws.Range("N7").Value = -275.77778
$ws.Range("H55").Value = 9143
$ws.Range("I55").Value = 5191
$ws.Range("J55").Value = 10625
$ws.Range("K55").Value = 5191
$ws.Range("L55").Value = 10625
$ws.Range("M55").Value = -4876
$ws.Range("N55").Value = -11255
$ws.Range("H99").Value = 2603.2964
$ws.Range("I99").Value = 1999.3572
$ws.Range("J99").Value = 3253.6924
$ws.Range("K99").Value = 1999.3572
$ws.Range("L99").Value = 3253.6924
$ws.Range("M99").Value = -501.3571999999999
$ws.Range("N99").Value = -6249.6924
$ws.Range("H122").Value = 1059
$ws.Range("I122").Value = 1157.7142
$ws.Range("J122").Value = 920.8
$ws.Range("K122").Value = 3473.1426
$ws.Range("L122").Value = 2762.4
$ws.Range("M122").Value = -1023.1426
$ws.Range("N122").Value = -7662.4
$ws.Range("H126").Value = 2603.2964
$ws.Range("I126").Value = 1999.3572
$ws.Range("J126").Value = 3253.6924
$ws.Range("K126").Value = 5998.071599999999
$ws.Range("L126").Value = 9761.0772
$ws.Range("M126").Value = -3528.071599999999
$ws.Range("N126").Value = -14701.0772
$ws.Range("H134").Value = 3472.5945
$ws.Range("J134").Value = 5879.7856
$ws.Range("L134").Value = 17639.3568
$ws.Range("N134").Value = -22709.3568

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 484.44446
$ws.Range("I5").Value = 340.0625
$ws.Range("J5").Value = 599.95
$ws.Range("K5").Value = 1020.1875
$ws.Range("L5").Value = 1799.85
$ws.Range("M5").Value = -908.1875
$ws.Range("N5").Value = -2023.85
$ws.Range("H43").Value = 7590
$ws.Range("J43").Value = 7590
$ws.Range("L43").Value = 22770
$ws.Range("N43").Value = -22998
$ws.Range("H44").Value = 362.25
$ws.Range("I44").Value = 239.6
$ws.Range("J44").Value = 566.6667
$ws.Range("K44").Value = 718.8
$ws.Range("L44").Value = 1700.0001
$ws.Range("M44").Value = -320.8
$ws.Range("N44").Value = -2496.0001
$ws.Range("H100").Value = 3112.5
$ws.Range("J100").Value = 4000
$ws.Range("L100").Value = 12000
$ws.Range("N100").Value = -13622
$ws.Range("H135").Value = 484.44446
$ws.Range("I135").Value = 340.0625
$ws.Range("J135").Value = 599.95
$ws.Range("K135").Value = 3060.5625
$ws.Range("L135").Value = 5399.55
$ws.Range("M135").Value = -525.5625
$ws.Range("N135").Value = -10469.55
$ws.Range("H140").Value = 5560833
$ws.Range("I140").Value = 8336832.5
$ws.Range("J140").Value = 8833.333000000001
$ws.Range("K140").Value = 25010497.5
$ws.Range("L140").Value = 26499.999
$ws.Range("M140").Value = -25005317.5
$ws.Range("N140").Value = -36859.999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2539.9644
$ws.Range("I122").Value = 2604.6667
$ws.Range("J122").Value = 2345.8572
$ws.Range("K122").Value = 7814.000100000001
$ws.Range("L122").Value = 7037.571599999999
$ws.Range("M122").Value = -5364.000100000001
$ws.Range("N122").Value = -11937.5716

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1031.7646
$ws.Range("I16").Value = 1036.5
$ws.Range("J16").Value = 1020.4
$ws.Range("K16").Value = 1036.5
$ws.Range("L16").Value = 1020.4
$ws.Range("M16").Value = -866.5
$ws.Range("N16").Value = -1360.4
$ws.Range("H22").Value = 304101.38
$ws.Range("I22").Value = 401077.7
$ws.Range("J22").Value = 1050.375
$ws.Range("K22").Value = 401077.7
$ws.Range("L22").Value = 1050.375
$ws.Range("M22").Value = -400782.7
$ws.Range("N22").Value = -1640.375
$ws.Range("H27").Value = 304101.38
$ws.Range("I27").Value = 401077.7
$ws.Range("J27").Value = 1050.375
$ws.Range("K27").Value = 401077.7
$ws.Range("L27").Value = 1050.375
$ws.Range("M27").Value = -400970.7
$ws.Range("N27").Value = -1264.375
$ws.Range("H29").Value = 2958
$ws.Range("I29").Value = 2958
$ws.Range("K29").Value = 2958
$ws.Range("M29").Value = -2663
$ws.Range("H46").Value = 62501436
$ws.Range("I46").Value = 76924290
$ws.Range("J46").Value = 2426.3333
$ws.Range("K46").Value = 76924290
$ws.Range("L46").Value = 2426.3333
$ws.Range("M46").Value = -76924102
$ws.Range("N46").Value = -2802.3333
$ws.Range("H132").Value = 8104.294
$ws.Range("I132").Value = 9365.556
$ws.Range("J132").Value = 6685.375
$ws.Range("K132").Value = 28096.668
$ws.Range("L132").Value = 20056.125
$ws.Range("M132").Value = -25566.668
$ws.Range("N132").Value = -25116.125
